$d = $word.ActiveDocument

# Update the title date line
$d.Content.Find.Execute("2025-10-09 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-10-10 Friday", 2)

# Update the division-problem table (5 columns, problem rows are 1, 5, 9, 13, 17)
$t = $d.Tables.Item(1)

$newValues = @{
    1  = @("50÷5=10, 0", "43÷8=5, 3", "36÷2=18, 0", "51÷6=8, 3", "78÷3=26, 0")
    5  = @("96÷9=10, 6", "90÷6=15, 0", "11÷7=1, 4", "25÷6=4, 1", "52÷6=8, 4")
    9  = @("13÷7=1, 6", "93÷7=13, 2", "55÷7=7, 6", "15÷4=3, 3", "84÷5=16, 4")
    13 = @("43÷9=4, 7", "29÷2=14, 1", "59÷6=9, 5", "48÷3=16, 0", "17÷2=8, 1")
    17 = @("89÷4=22, 1", "85÷8=10, 5", "12÷8=1, 4", "84÷8=10, 4", "36÷5=7, 1")
}

foreach ($rowIndex in $newValues.Keys) {
    $cols = $newValues[$rowIndex]
    for ($col = 1; $col -le 5; $col++) {
        $t.Cell($rowIndex, $col).Range.Text = $cols[$col - 1]
    }
}
